$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 177, pushing existing rows 177-238 down to 178-239.
$ws.Rows.Item(177).Insert()

# Populate the newly inserted row 177 with the new weekly price record.
$ws.Range("A177").Value = 8
$ws.Range("B177").Value = "Terminal La Palmera de La Serena"
$ws.Range("C177").Value = "Coquimbo"
$ws.Range("D177").Value = 44468
$ws.Range("E177").Value = 4
$ws.Range("F177").Value = 100114001
$ws.Range("G177").Value = "Papa"
$ws.Range("H177").Value = "Cardinal"
$ws.Range("I177").Value = "1a (cosecha)"
$ws.Range("J177").Value = 2900
$ws.Range("K177").Value = 12000
$ws.Range("L177").Value = 13000
$ws.Range("M177").Value = 12500
$ws.Range("N177").Value = "$/saco 25 kilos"
$ws.Range("O177").Value = "Provincia del Elquí"
$ws.Range("P177").Value = 500
$ws.Range("Q177").Value = 25
$ws.Range("R177").Value = "Hortaliza"
